{"js": "// Apply the per-cell text replacements described by the diff:\n// the date line and each of the 26 two-digit multiplication problems\n// are replaced with new values, in document order.\nconst replacements = [\n  [\"2024-05-04 Saturday\", \"2024-05-05 Sunday\"],\n  [\"62\u00d761=\", \"93\u00d775=\"],\n  [\"54\u00d768=\", \"12\u00d739=\"],\n  [\"64\u00d735=\", \"72\u00d768=\"],\n  [\"70\u00d756=\", \"87\u00d715=\"],\n  [\"29\u00d730=\", \"20\u00d796=\"],\n  [\"98\u00d766=\", \"54\u00d712=\"],\n  [\"45\u00d741=\", \"78\u00d796=\"],\n  [\"41\u00d798=\", \"21\u00d797=\"],\n  [\"22\u00d793=\", \"75\u00d764=\"],\n  [\"81\u00d779=\", \"52\u00d786=\"],\n  [\"71\u00d763=\", \"43\u00d787=\"],\n  [\"99\u00d719=\", \"44\u00d767=\"],\n  [\"25\u00d717=\", \"21\u00d796=\"],\n  [\"83\u00d771=\", \"54\u00d741=\"],\n  [\"44\u00d785=\", \"62\u00d798=\"],\n  [\"25\u00d734=\", \"15\u00d747=\"],\n  [\"90\u00d796=\", \"81\u00d795=\"],\n  [\"47\u00d718=\", \"30\u00d727=\"],\n  [\"43\u00d760=\", \"66\u00d766=\"],\n  [\"77\u00d760=\", \"53\u00d783=\"],\n  [\"53\u00d712=\", \"13\u00d782=\"],\n  [\"74\u00d793=\", \"66\u00d799=\"],\n  [\"98\u00d782=\", \"88\u00d766=\"],\n  [\"72\u00d798=\", \"99\u00d717=\"],\n  [\"80\u00d798=\", \"42\u00d727=\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the per-cell text replacements described by the diff:\n# the date line and each of the 26 two-digit multiplication problems\n# are replaced with new values, in document order.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2024-05-04 Saturday\", \"2024-05-05 Sunday\"),\n    @(\"62\u00d761=\", \"93\u00d775=\"),\n    @(\"54\u00d768=\", \"12\u00d739=\"),\n    @(\"64\u00d735=\", \"72\u00d768=\"),\n    @(\"70\u00d756=\", \"87\u00d715=\"),\n    @(\"29\u00d730=\", \"20\u00d796=\"),\n    @(\"98\u00d766=\", \"54\u00d712=\"),\n    @(\"45\u00d741=\", \"78\u00d796=\"),\n    @(\"41\u00d798=\", \"21\u00d797=\"),\n    @(\"22\u00d793=\", \"75\u00d764=\"),\n    @(\"81\u00d779=\", \"52\u00d786=\"),\n    @(\"71\u00d763=\", \"43\u00d787=\"),\n    @(\"99\u00d719=\", \"44\u00d767=\"),\n    @(\"25\u00d717=\", \"21\u00d796=\"),\n    @(\"83\u00d771=\", \"54\u00d741=\"),\n    @(\"44\u00d785=\", \"62\u00d798=\"),\n    @(\"25\u00d734=\", \"15\u00d747=\"),\n    @(\"90\u00d796=\", \"81\u00d795=\"),\n    @(\"47\u00d718=\", \"30\u00d727=\"),\n    @(\"43\u00d760=\", \"66\u00d766=\"),\n    @(\"77\u00d760=\", \"53\u00d783=\"),\n    @(\"53\u00d712=\", \"13\u00d782=\"),\n    @(\"74\u00d793=\", \"66\u00d799=\"),\n    @(\"98\u00d782=\", \"88\u00d766=\"),\n    @(\"72\u00d798=\", \"99\u00d717=\"),\n    @(\"80\u00d798=\", \"42\u00d727=\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2)\n}\n"}
